# Update the cryptos price/volume table (rows 2-51) to the refreshed
# scrape values, matching the upstream GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.004.42'
$ws.Range("E2").Value = '  -1.29%  '

$ws.Range("D3").Value = '2.167.55'
$ws.Range("E3").Value = '  -2.46%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.19'
$ws.Range("E5").Value = '  -1.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("E6").Value = '  -2.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.84'
$ws.Range("E7").Value = '  -6.85%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  -3.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.70'
$ws.Range("E10").Value = '  +2.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0922'
$ws.Range("E11").Value = '  -4.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.17'
$ws.Range("E12").Value = '  -14.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  -1.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.79'
$ws.Range("E14").Value = '  -2.19%  '

$ws.Range("D15").Value = '2.496.05'
$ws.Range("E15").Value = '  -2.10%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.19'
$ws.Range("E16").Value = '  -4.38%  '

$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.847'
$ws.Range("E17").Value = '  -0.84%  '

$ws.Range("D18").Value = '2.172.61'
$ws.Range("E18").Value = '  -1.73%  '

$ws.Range("D19").Value = '40.870.90'
$ws.Range("E19").Value = '  -1.52%  '

$ws.Range("D20").Value = '0.0₃0935'
$ws.Range("E20").Value = '  -3.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.24'
$ws.Range("E21").Value = '  -1.71%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("E22").Value = '  -2.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.95'
$ws.Range("E23").Value = '  -2.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").Value = '  -7.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("B26").Value = 'WEMIXToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.70'
$ws.Range("E26").Value = '  -5.03%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.13'
$ws.Range("E27").Value = '  +6.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.41'
$ws.Range("E28").Value = '  -3.12%  '

$ws.Range("E29").Value = '  -5.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.10'
$ws.Range("E30").Value = '  -1.81%  '

$ws.Range("E31").Value = '  -8.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.11'
$ws.Range("E32").Value = '  -2.27%  '

$ws.Range("E33").Value = '  +1.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.62'
$ws.Range("E34").Value = '  +1.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0740'
$ws.Range("E35").Value = '  +3.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.120'
$ws.Range("E36").Value = '  -3.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.53'
$ws.Range("E37").Value = '  -3.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.94'
$ws.Range("E38").Value = '  +0.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.16'
$ws.Range("E39").Value = '  -6.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0299'
$ws.Range("E40").Value = '  +2.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.17'
$ws.Range("E41").Value = '  -4.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.44'
$ws.Range("E42").Value = '  -8.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.82'
$ws.Range("E43").Value = '  -0.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '59.94'
$ws.Range("E44").Value = '  -12.67%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.99'
$ws.Range("E45").Value = '  -6.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.190'
$ws.Range("E46").Value = '  -8.71%  '

$ws.Range("B47").Value = 'BinanceUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.05%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.43'
$ws.Range("E48").Value = '  -3.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0984'
$ws.Range("E49").Value = '  -2.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.14'
$ws.Range("E50").Value = '  -1.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.14'
$ws.Range("E51").Value = '  -3.77%  '
